$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new columns at D:E, shifting existing D:K data to F:M.
$ws.Columns("D:E").Insert()

# 2) Copy the number/date formatting from the (now shifted) old "D" column,
#    which now lives in column F, into the two new D:E columns so the new
#    cells carry the same styling (date format for header rows, number
#    format for data rows) as their neighbours.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate the two new columns (D = newest quarter, E = next quarter)
#    with the values for each row.
$deData = @{
    7 = @(43463, 43372)
    8 = @(299000, 296800)
    9 = @(49500, 51900)
    10 = @(249500, 244900)
    11 = @($null, $null)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(250500, 251300)
    18 = @(48500, 45500)
    19 = @($null, $null)
    20 = @(900, 1000)
    21 = @(53600, 50600)
    22 = @(0, 0)
    23 = @(49500, 46500)
    24 = @(17100, 15500)
    25 = @(0, 0)
    26 = @(32300, 31000)
    27 = @(32300, 31000)
    28 = @(0, 0)
    29 = @("NA", "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-900, -1000)
    33 = @(32300, 31000)
    34 = @(0, 0)
    35 = @(32300, 31000)
    38 = @(43463, 43372)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(214300, 239800)
    42 = @(63500, 81700)
    43 = @(0, 0)
    44 = @(81900, 81100)
    45 = @(32500, 32200)
    46 = @(392300, 434800)
    47 = @(0, 0)
    48 = @(92000, 93600)
    49 = @(48600, 49000)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(21500, 22000)
    53 = @(0, 0)
    54 = @(554500, 599400)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(9900, 14500)
    58 = @(0, 0)
    59 = @(138700, 138400)
    60 = @(148700, 152900)
    61 = @(0, 0)
    62 = @(14600, 10900)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(163300, 163800)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(329500, 362300)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(391100, 435600)
    77 = @(0, 0)
    80 = @(43463, 43372)
    81 = @(32300, 31000)
    82 = @($null, $null)
    83 = @(4100, 4000)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(40400, 54500)
    90 = @($null, $null)
    91 = @(-2600, -2200)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(15600, -41400)
    95 = @($null, $null)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-80200, -22700)
    101 = @(-1400, -7100)
    102 = @(-25500, -16700)
}

foreach ($r in $deData.Keys) {
    $vals = $deData[$r]
    if ($null -ne $vals[0]) {
        $ws.Cells.Item($r, 4).Value = $vals[0]
    }
    if ($null -ne $vals[1]) {
        $ws.Cells.Item($r, 5).Value = $vals[1]
    }
}

# 4) A couple of rows also received corrected historical figures (not just
#    a straight shift) further along the row; apply those explicitly.
# Row 58: columns F:J (previously "NA") now contain 0.
$ws.Range("F58:J58").Value = 0
# Row 91: columns I:J were restated.
$ws.Cells.Item(91, 9).Value = -2200
$ws.Cells.Item(91, 10).Value = -3600

# 5) Approximate the new column widths for the two inserted columns so
#    they resemble their neighbours.
$ws.Columns("D:E").ColumnWidth = 14.6640625

Write-Host "Edit applied"
